$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename VenueLatitude/VenueLongitude headers to Latitude/Longitude,
# and add a new "Visited" header in column T
$ws.Range("Q1").Value = "Latitude"
$ws.Range("R1").Value = "Longitude"
$ws.Range("T1").Value = "Visited"

# Populate the new Visited column with boolean flags for each ballpark row
$ws.Range("T2").Value = $false
$ws.Range("T3").Value = $false
$ws.Range("T4").Value = $false
$ws.Range("T5").Value = $false
$ws.Range("T6").Value = $false
$ws.Range("T7").Value = $false
$ws.Range("T8").Value = $false
$ws.Range("T9").Value = $false
$ws.Range("T10").Value = $false
$ws.Range("T11").Value = $false
$ws.Range("T12").Value = $false
$ws.Range("T13").Value = $false
$ws.Range("T14").Value = $false
$ws.Range("T15").Value = $false
$ws.Range("T16").Value = $false
$ws.Range("T17").Value = $false
$ws.Range("T18").Value = $false
$ws.Range("T19").Value = $false
$ws.Range("T20").Value = $false
$ws.Range("T21").Value = $false
$ws.Range("T22").Value = $false
$ws.Range("T23").Value = $false
$ws.Range("T24").Value = $true
$ws.Range("T25").Value = $false
$ws.Range("T26").Value = $false
$ws.Range("T27").Value = $true
$ws.Range("T28").Value = $false
$ws.Range("T29").Value = $false
$ws.Range("T30").Value = $true
$ws.Range("T31").Value = $true
$ws.Range("T32").Value = $true
$ws.Range("T33").Value = $true
$ws.Range("T34").Value = $true
$ws.Range("T35").Value = $true
$ws.Range("T36").Value = $true
$ws.Range("T37").Value = $true
$ws.Range("T38").Value = $true
$ws.Range("T39").Value = $true
$ws.Range("T40").Value = $true
$ws.Range("T41").Value = $true
$ws.Range("T42").Value = $true
$ws.Range("T43").Value = $true
$ws.Range("T44").Value = $true
$ws.Range("T45").Value = $true
$ws.Range("T46").Value = $true
$ws.Range("T47").Value = $true
$ws.Range("T48").Value = $true
$ws.Range("T49").Value = $true
$ws.Range("T50").Value = $true
$ws.Range("T51").Value = $true
$ws.Range("T52").Value = $true
$ws.Range("T53").Value = $true
$ws.Range("T54").Value = $true
$ws.Range("T55").Value = $true
$ws.Range("T56").Value = $true
$ws.Range("T57").Value = $true
$ws.Range("T58").Value = $true
$ws.Range("T59").Value = $true
$ws.Range("T60").Value = $true
$ws.Range("T61").Value = $true

# Restore the selection/scroll state shown in the saved workbook
$ws.Range("S22").Select() | Out-Null
